# Insert a new slide ("Numbers can also be assigned like this") right
# before the final "Any Questions?" slide (slide index 8 of 8).
# Using layout 2 ("Title and Content") to match every other slide in
# this deck.
$p = $ppt.ActivePresentation

$newSlide = $p.Slides.Add(8, 2)

# --- Title placeholder -----------------------------------------------
$title = $newSlide.Shapes.Item(1).TextFrame.TextRange
$title.Text = "Numbers can also be assigned like this "
$title.Font.Bold = $true

# --- Body / content placeholder ---------------------------------------
$body = $newSlide.Shapes.Item(2).TextFrame.TextRange

# Paragraph 1: "-  a,b,c,d = 3,4,5,6"
$body.Text = [char]0x2013 + "  "
$body.Font.Size = 36
$body.ParagraphFormat.Bullet.Visible = $false

$run = $body.InsertAfter("a,b,c,d")
$run.Font.Size = 36

$run = $run.InsertAfter(" = 3,4,5,6")
$run.Font.Size = 36

# Paragraph 2: "    print (a,b,c,d)"
$run = $run.InsertAfter("`r    print (")
$run.Font.Size = 36
$run.ParagraphFormat.Bullet.Visible = $false

$run = $run.InsertAfter("a,b,c,d")
$run.Font.Size = 36

$run = $run.InsertAfter(")")
$run.Font.Size = 36

# Paragraph 3: blank line
$run = $run.InsertAfter("`r")
$run.ParagraphFormat.Bullet.Visible = $false

# Paragraph 4: blank trailing paragraph
$run = $run.InsertAfter("`r")
$run.ParagraphFormat.Bullet.Visible = $false
